# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets.
# Values derived from the canonical OOXML diff for this commit.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1614
$ws.Range("G2").Value = 45
$ws.Range("F3").Value = 849
$ws.Range("F4").Value = 257
$ws.Range("F5").Value = 74
$ws.Range("F6").Value = 1165
$ws.Range("F7").Value = 773
$ws.Range("F8").Value = 810
$ws.Range("F9").Value = 1486
$ws.Range("F11").Value = 1046
$ws.Range("F12").Value = 31
$ws.Range("F14").Value = 196
$ws.Range("F15").Value = 54
$ws.Range("F16").Value = 495
$ws.Range("F17").Value = 46
$ws.Range("F19").Value = 6
$ws.Range("F22").Value = 566
$ws.Range("F24").Value = 35
$ws.Range("F27").Value = 254
$ws.Range("F28").Value = 188

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1016
$ws.Range("F5").Value = 275
$ws.Range("F7").Value = 149
$ws.Range("F9").Value = 591
$ws.Range("F10").Value = 87

# ---- Sheet "本地生活" ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 261

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 261
$ws.Range("F3").Value = 1614
$ws.Range("G3").Value = 45
$ws.Range("F5").Value = 849
$ws.Range("F6").Value = 257
$ws.Range("F7").Value = 1016
$ws.Range("F8").Value = 74
$ws.Range("F9").Value = 1165
$ws.Range("F10").Value = 773
$ws.Range("F11").Value = 810
$ws.Range("F12").Value = 1486
$ws.Range("F14").Value = 1046
$ws.Range("F15").Value = 31
$ws.Range("F17").Value = 196
$ws.Range("F18").Value = 54
$ws.Range("F19").Value = 495
$ws.Range("F20").Value = 46
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 275
$ws.Range("F28").Value = 149
$ws.Range("F29").Value = 149
$ws.Range("F30").Value = 566
$ws.Range("F32").Value = 35
$ws.Range("F35").Value = 254
$ws.Range("F37").Value = 188
$ws.Range("F38").Value = 591
$ws.Range("F39").Value = 87
$ws.Range("F40").Value = 87
